$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BBL")

$ws.Range("D8").Value = 43129000

$ws.Range("D9").Value = 7867000
$ws.Range("E9").Value = 6224000
$ws.Range("F9").Value = 6615000
$ws.Range("G9").Value = 7679000
$ws.Range("H9").Value = 9520000
$ws.Range("I9").Value = 8736000
$ws.Range("J9").Value = 26682000

$ws.Range("D10").Value = 35262000
$ws.Range("E10").Value = 29911000
$ws.Range("F10").Value = 21952000
$ws.Range("G10").Value = 36957000
$ws.Range("H10").Value = 47242000
$ws.Range("I10").Value = 45124000
$ws.Range("J10").Value = 43795000

$ws.Range("D17").Value = 27133000

$ws.Range("D20").Value = -157000

$ws.Range("D21").Value = 22127000

$ws.Range("D22").Value = 1088000

$ws.Range("D32").Value = 157000

$ws.Range("G91").Value = -12763000
$ws.Range("H91").Value = -16210000
$ws.Range("I91").Value = -22425000
$ws.Range("J91").Value = -21130000
